# Regenerate localization-status report: update translation status for the
# two files that moved from "Ready for handoff" into "In Translation"
# (1b34b1d9-c08a-4ed6-b90b-4612acc7419f.md and
#  1c6ff0da-558e-47d1-8404-6c35d564d400.md) across the Overview summary
# sheet as well as the per-locale (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) and de-de (col F) status columns ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus   # 1b34b1d9...md / zh-cn
$overview.Range("F3").Value = $newStatus   # 1b34b1d9...md / de-de
$overview.Range("E4").Value = $newStatus   # 1c6ff0da...md / zh-cn
$overview.Range("F4").Value = $newStatus   # 1c6ff0da...md / de-de

# --- zh-cn detail sheet: Status column (col C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus       # 1b34b1d9...md
$zhcn.Range("C4").Value = $newStatus       # 1c6ff0da...md

# --- de-de detail sheet: Status column (col C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus       # 1b34b1d9...md
$dede.Range("C4").Value = $newStatus       # 1c6ff0da...md

Write-Host "Updated status to 'In Translation' for 1b34b1d9 and 1c6ff0da rows on Overview, zh-cn, de-de sheets."
